$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "iaest-measure:estrato"
$ws.Range("G2").Value = "iaest-measure:direccion-provincial-nombre"

$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"

$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"

$ws.Range("F5").ClearContents()
